$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.242.89"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -15.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.252.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -22.41%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "429.73"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -18.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "114.56"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -20.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.449"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -17.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.255.17"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -22.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.07"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -16.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0831"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -22.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.292"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -18.44%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -7.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.616.44"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -23.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "51.258.00"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -15.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.06"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -19.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000111"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -20.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.246.26"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -22.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.83"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -22.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "286.99"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -18.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.990"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.33"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -28.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.87"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -25.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "52.18"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -19.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.356"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -21.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.325.27"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -23.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.134"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -23.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.996"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.54"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -16.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "141.95"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0608"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -28.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "16.17"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -17.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.28"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -23.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.50"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -18.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.995"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -25.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "31.39"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -16.34%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.941"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -21.19%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.721"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -27.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.12"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.544"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -16.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -18.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0478"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -17.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.829.10"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -19.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.09"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -25.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0196"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -17.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0778"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -15.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.62"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.43%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.82"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -22.55%  "
